$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.870848894119263
$ws.Range("B1").Value = 2.003752470016479
$ws.Range("C1").Value = 2.100384473800659
$ws.Range("D1").Value = 2.824164390563965
$ws.Range("E1").Value = 3.356984853744507
